# Atualizando a minha experiência
#
# Adds a new bullet item ("Desenvolvimento Web") to the end of the
# "Habilidades" (Skills) list, right after the existing
# "Back-End - intermediário" bullet, as the new last paragraph of the
# document body.

$d = $word.ActiveDocument

# Avoid double-applying the edit if the script were ever run twice.
$already = $d.Content.Find.Execute("Desenvolvimento Web", $true, $false, $false, $false, $false,
                                    $true, 1, $false, "", 0)

if (-not $already) {
    $rng = $d.Content
    $found = $rng.Find.Execute("Back-End - intermediário ", $true, $false, $false, $false, $false,
                                $true, 1, $false, "", 0)

    if (-not $found) {
        # Fall back to a slightly looser search (without the trailing space)
        # in case whitespace normalization differs.
        $rng = $d.Content
        $found = $rng.Find.Execute("Back-End - intermediário", $true, $false, $false, $false, $false,
                                    $true, 1, $false, "", 0)
    }

    if ($found) {
        # Collapse to right after the matched text (still before the very
        # last paragraph mark of the document) and add a new paragraph.
        # The new paragraph naturally inherits the same paragraph style
        # ("Parágrafo da Lista") and list numbering (numId 5) from the
        # paragraph it follows.
        $rng.Collapse(0)
        $rng.InsertAfter([char]13 + "Desenvolvimento Web")
    } else {
        # Last-resort fallback: append at the very end of the document.
        $endRng = $d.Range($d.Content.End, $d.Content.End)
        $endRng.InsertAfter([char]13 + "Desenvolvimento Web")
    }
}

$d.Saved = $false
